# "update to v3 change author order"
# On the Personnel sheet, the author/creator rows for Tatiana Rynearson
# (row 3) and Diana Fontaine (row 4) are swapped - i.e. Diana's record now
# appears in row 3 and Tatiana's in row 4 (columns A:F, which carry the
# per-person fields; the shared G:J columns - role/org/funding - are left
# untouched since both rows already held identical values there). The E-column
# mailto hyperlinks follow the person they belong to, and the row heights
# (which differ because row 3 originally carried extra top padding) move
# along with the content. Finally, the previously-active sheet/tab
# (ColumnHeadersIntegrated) is deactivated and the Personnel sheet - with its
# newly-edited row 3 selected - becomes the active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# --- swap the per-person fields (A:F) between row 3 and row 4 ---
$row3 = @()
$row4 = @()
for ($c = 1; $c -le 6; $c++) {
    $row3 += $ws.Cells.Item(3, $c).Value2
    $row4 += $ws.Cells.Item(4, $c).Value2
}
for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item(3, $c).Value = $row4[$c - 1]
    $ws.Cells.Item(4, $c).Value = $row3[$c - 1]
}

# the taller row (ht 15.6) moves from row 3 to row 4 along with the content
$ws.Rows.Item(3).RowHeight = 14.4
$ws.Rows.Item(4).RowHeight = 15.6

# --- the E3/E4 mailto hyperlinks follow their person to the new row ---
$hyperlinkE3 = $null
$hyperlinkE4 = $null
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$E$3') { $hyperlinkE3 = $h }
    if ($addr -eq '$E$4') { $hyperlinkE4 = $h }
}
$addressAtE3 = $hyperlinkE3.Address
$addressAtE4 = $hyperlinkE4.Address
$hyperlinkE3.Address = $addressAtE4
$hyperlinkE4.Address = $addressAtE3

# --- Personnel becomes the active sheet, with the edited row selected ---
$ws.Activate()
$ws.Rows.Item(3).Select()
